$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new rows before the existing row 2, shifting old data down to rows 11-29
$ws.Range("A2:C10").EntireRow.Insert()
$ws.Range("A2:C10").ClearFormats()

# Populate the newly inserted rows (2-10) with the new accelerometer samples
$ws.Cells.Item(2, 1).Value = -0.4484343528747558
$ws.Cells.Item(2, 2).Value = 0.5805449485778809
$ws.Cells.Item(2, 3).Value = -0.5587977170944214
$ws.Cells.Item(3, 1).Value = -0.6014323234558105
$ws.Cells.Item(3, 2).Value = 0.5997557640075684
$ws.Cells.Item(3, 3).Value = -0.4316743612289428
$ws.Cells.Item(4, 1).Value = -0.2612781524658203
$ws.Cells.Item(4, 2).Value = 0.5193090438842773
$ws.Cells.Item(4, 3).Value = -0.4944255352020263
$ws.Cells.Item(5, 1).Value = -0.2169137001037597
$ws.Cells.Item(5, 2).Value = 0.3676133155822754
$ws.Cells.Item(5, 3).Value = -0.733814001083374
$ws.Cells.Item(6, 1).Value = -0.3025293350219726
$ws.Cells.Item(6, 2).Value = 0.4683008193969726
$ws.Cells.Item(6, 3).Value = -0.587003767490387
$ws.Cells.Item(7, 1).Value = -0.4149298667907715
$ws.Cells.Item(7, 2).Value = 0.4677276611328125
$ws.Cells.Item(7, 3).Value = -0.731619656085968
$ws.Cells.Item(8, 1).Value = -0.5511326789855957
$ws.Cells.Item(8, 2).Value = 0.6498098373413086
$ws.Cells.Item(8, 3).Value = -0.522668182849884
$ws.Cells.Item(9, 1).Value = -0.5137066841125488
$ws.Cells.Item(9, 2).Value = 0.4998054504394531
$ws.Cells.Item(9, 3).Value = -0.6402766704559326
$ws.Cells.Item(10, 1).Value = -0.6838326454162598
$ws.Cells.Item(10, 2).Value = 0.6059346199035645
$ws.Cells.Item(10, 3).Value = -0.2089821100234985

# Append one more row (31) of new data at the end of the table
$ws.Cells.Item(31, 1).Value = -0.143467903137207
$ws.Cells.Item(31, 2).Value = 0.759878396987915
$ws.Cells.Item(31, 3).Value = -0.6699965000152588
